$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 48
$ws1.Range("F9").Value = 1241
$ws1.Range("F11").Value = 288
$ws1.Range("F12").Value = 1088
$ws1.Range("F14").Value = 6861
$ws1.Range("F18").Value = 7765
$ws1.Range("F20").Value = 43
$ws1.Range("F21").Value = 4629
$ws1.Range("F23").Value = 2261
$ws1.Range("F24").Value = 956
$ws1.Range("F26").Value = 241
$ws1.Range("F27").Value = 361
$ws1.Range("F30").Value = 266
$ws1.Range("F31").Value = 228
$ws1.Range("F33").Value = 1966
$ws1.Range("F34").Value = 15
$ws1.Range("F35").Value = 221
$ws1.Range("F37").Value = 525
$ws1.Range("F39").Value = 1345
$ws1.Range("F40").Value = 16
$ws1.Range("F41").Value = 2054
$ws1.Range("F42").Value = 2174

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 21
$ws2.Range("F5").Value = 11

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1256
$ws3.Range("F4").Value = 82

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 48
$ws4.Range("F4").Value = 1256
$ws4.Range("F5").Value = 82
$ws4.Range("F11").Value = 1241
$ws4.Range("F13").Value = 288
$ws4.Range("F14").Value = 1088
$ws4.Range("F16").Value = 6861
$ws4.Range("F20").Value = 7765
$ws4.Range("F22").Value = 43
$ws4.Range("F23").Value = 4629
$ws4.Range("F25").Value = 2261
$ws4.Range("F26").Value = 956
$ws4.Range("F28").Value = 241
$ws4.Range("F29").Value = 361
$ws4.Range("F33").Value = 21
$ws4.Range("F34").Value = 266
$ws4.Range("F36").Value = 1966
$ws4.Range("F37").Value = 15
$ws4.Range("F38").Value = 221
$ws4.Range("F40").Value = 525
$ws4.Range("F42").Value = 11
$ws4.Range("F43").Value = 1345
$ws4.Range("F44").Value = 16
$ws4.Range("F45").Value = 2054
$ws4.Range("F47").Value = 2174
